$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Youtube names scrapping"

# Shift the existing "video titles" column from A to B by inserting a new column A
$ws.Columns("A:A").Insert()

# Apply the header style (copied from the now-shifted B1 header cell) down A2:A25
$ws.Range("B1").Copy()
$ws.Range("A2:A25").PasteSpecial(-4122)

# Fill in column A with the new numeric index values
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 3
$ws.Range("A4").Value = 6
$ws.Range("A5").Value = 9
$ws.Range("A6").Value = 12
$ws.Range("A7").Value = 15
$ws.Range("A8").Value = 18
$ws.Range("A9").Value = 20
$ws.Range("A10").Value = 23
$ws.Range("A11").Value = 26
$ws.Range("A12").Value = 29
$ws.Range("A13").Value = 32
$ws.Range("A14").Value = 35
$ws.Range("A15").Value = 38
$ws.Range("A16").Value = 41
$ws.Range("A17").Value = 44
$ws.Range("A18").Value = 47
$ws.Range("A19").Value = 50
$ws.Range("A20").Value = 53
$ws.Range("A21").Value = 56
$ws.Range("A22").Value = 58
$ws.Range("A23").Value = 61
$ws.Range("A24").Value = 64
$ws.Range("A25").Value = 67

# Replace column B (the shifted "video titles") rows 2-25 with the new video titles
$ws.Range("B2").Value = "Web Scraping Facebook with Selenium - AUTOMATED BOT"
$ws.Range("B3").Value = "Scrape Any Facebook Group's Posts with Selenium & BeautifulSoup (Free, works for private groups!)"
$ws.Range("B4").Value = "Web Scraping using Python and Selenium | Scrape Facebook | Part 5 | Data Making | DM | DataMaking"
$ws.Range("B5").Value = "How to SCRAPE DYNAMIC websites with Selenium"
$ws.Range("B6").Value = "Python Webcrape Facebook Marketplace"
$ws.Range("B7").Value = "Python Selenium- Facebook Marketplace Webscrape Part2"
$ws.Range("B8").Value = "Facebook Post Comments Scraper Using Python Selenium"
$ws.Range("B9").Value = "Top 3 FREE Methods using Email Extractor for FREE Email Marketing email extractor free"
$ws.Range("B10").Value = "Facebook Page Scraping | Scraping B2B Pages For Emails And Phone Numbers In 2021 💲 Scrapebox 👈"
$ws.Range("B11").Value = "Scrape Emails From Facebook Business Pages With Scrapebox : Updated Method For 2020"
$ws.Range("B12").Value = "Scrape Any Website Without Code | Generate Leads | Collect Any Data"
$ws.Range("B13").Value = "Web Scraping Instagram with Selenium"
$ws.Range("B14").Value = "Python Selenium Tutorial #1 - Web Scraping, Bots & Testing"
$ws.Range("B15").Value = "Beautifulsoup vs Selenium vs Scrapy - Which tool for web scraping in 2021?"
$ws.Range("B16").Value = "(Python)Get Facebook Page Data using Python Scraping (Selenium)"
$ws.Range("B17").Value = "Python Web Scraping - Should I use Selenium, Beautiful Soup or Scrapy? [2020]"
$ws.Range("B18").Value = "How to scrape INFINITE scrolling pages using Python and Selenium (2 Methods)"
$ws.Range("B19").Value = "👥 Facebook Data Scraping & Legal Deep Dive - Free HAR File Tool to Responsibly Extract Facebook Data"
$ws.Range("B20").Value = "#8 How To Automate Facebook Registration /Login Using Selenium Webdriver-Selenium Python Script"
$ws.Range("B21").Value = "Web Scraping EP4: bruteforcing facebook with selenium"
$ws.Range("B22").Value = "How To Scrape Facebook For 1000's Of Leads For Free. No-Code"
$ws.Range("B23").Value = "Facebook sues developer | Let's talk about data scraping"
$ws.Range("B24").Value = "How to Scrape Websites Without Getting Blacklisted or Blocked"
$ws.Range("B25").Value = "WEB SCRAPING made simple with JAVASCRIPT tutorial"
